$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the formatting
# used by the existing header cells (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add new data cells I2 and J2 for row 2
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
